$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

# Map of (row, column, newText) for each data cell that changed.
# The table has data only in rows 1, 5, 9, 13, 17 (1-indexed); other rows are blank spacer rows.
$updates = @(
    @{ Row = 1;  Col = 1; Text = "58÷4=14, 2" },
    @{ Row = 1;  Col = 2; Text = "51÷7=7, 2" },
    @{ Row = 1;  Col = 3; Text = "23÷7=3, 2" },
    @{ Row = 1;  Col = 4; Text = "94÷6=15, 4" },
    @{ Row = 1;  Col = 5; Text = "54÷8=6, 6" },

    @{ Row = 5;  Col = 1; Text = "23÷3=7, 2" },
    @{ Row = 5;  Col = 2; Text = "98÷6=16, 2" },
    @{ Row = 5;  Col = 3; Text = "18÷7=2, 4" },
    @{ Row = 5;  Col = 4; Text = "37÷3=12, 1" },
    @{ Row = 5;  Col = 5; Text = "65÷6=10, 5" },

    @{ Row = 9;  Col = 1; Text = "70÷8=8, 6" },
    @{ Row = 9;  Col = 2; Text = "88÷6=14, 4" },
    @{ Row = 9;  Col = 3; Text = "75÷6=12, 3" },
    @{ Row = 9;  Col = 4; Text = "30÷6=5, 0" },
    @{ Row = 9;  Col = 5; Text = "54÷7=7, 5" },

    @{ Row = 13; Col = 1; Text = "35÷5=7, 0" },
    @{ Row = 13; Col = 2; Text = "80÷9=8, 8" },
    @{ Row = 13; Col = 3; Text = "53÷8=6, 5" },
    @{ Row = 13; Col = 4; Text = "95÷9=10, 5" },
    @{ Row = 13; Col = 5; Text = "24÷7=3, 3" },

    @{ Row = 17; Col = 1; Text = "72÷4=18, 0" },
    @{ Row = 17; Col = 2; Text = "80÷7=11, 3" },
    @{ Row = 17; Col = 3; Text = "97÷3=32, 1" },
    @{ Row = 17; Col = 4; Text = "89÷3=29, 2" },
    @{ Row = 17; Col = 5; Text = "77÷2=38, 1" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.Text
}
